$d = $word.ActiveDocument

$questions = @(
    "Have you created your own algorithms before? How did you develop them?",
    "How do you approach data cleaning in Python?",
    "When working with a large dataset, how do you account for outliers? Missing values or transformation?",
    "What are some sorting algorithms you have used in the R language?",
    "If your data team uses Hadoop, how would you integrate it with R for enhanced data analysis?",
    "Describe your proficiency in SQL?",
    "What are some successful projects you have completed in SQL? What made them successful?",
    "What applications would you use the recommender system with?",
    "How do you apply univariate analysis?",
    "What process do you use to define the number of cluster values within a clustering algorithm?",
    "Explain your understanding of auto-encoders?",
    "How would you apply the batch normalisation process to organise and analyse data systems?",
    "What machine learning library do you feel is most beneficial for supervised learning projects?",
    "What steps do you take before you apply machine learning algorithms?",
    "How would you resolve unbalanced binary classifications?",
    "What are the benefits of box plots when visualization big data?",
    "Describe the regularisation method you can apply when implementing training data?",
    "How do you select metrics for cross-validation?",
    "How would you evaluate a predictive model from multiple regression analysis?",
    "When would you use random forests over a support vector machines?"
)

# Locate the paragraph ending "...critical data?" -- the last question in the
# existing list -- so we can append the new questions right after it.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*critical data?*") {
        $anchor = $p
    }
}

foreach ($q in $questions) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $q
}

# The paragraph immediately after the newly-inserted questions is the
# pre-existing blank "ListParagraph" paragraph; turn it into a plain
# (non-list) paragraph that just keeps a left indent, matching the target.
$blank = $anchor.Next()
$blank.Range.ParagraphFormat.Style = "Normal"
$blank.Range.ParagraphFormat.LeftIndent = 18
